# Update "想去人数" (want-to-go count) figures in column F for the
# data sheets "展览" (index 1) and "全部类型" (index 4). These two
# sheets mirror the same underlying data, so both need updating.
# Sheets "演出" (index 2) and "本地生活" (index 3) contain only header
# rows with no data and are left untouched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Row -> new F-column value, shared by both sheets (row 19 differs,
# handled separately below).
$rows  = @(2, 3, 4, 5, 9, 10, 12, 13, 14, 16)
$newVals = @(7413, 7374, 95, 183, 94, 128, 97, 668, 522, 30)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $v = $newVals[$i]
    $ws1.Cells.Item($r, 6).Value = $v
    $ws4.Cells.Item($r, 6).Value = $v
}

# F19 differs between the two sheets.
$ws1.Cells.Item(19, 6).Value = 75
$ws4.Cells.Item(19, 6).Value = 76
